# vk_photos: populate a small "file_name"/"size" table and style the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "file_name"
$ws.Range("B1").Value = "size"

# Data rows
$data = @(
  @("6_298279696.jpg", "z"),
  @("5_298279713.jpg", "z"),
  @("1_298279756.jpg", "z"),
  @("0_457239382.jpg", "z"),
  @("0_457239377.jpg", "z")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Style the header: bold, centered/top-aligned, thin box border all around.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108  # xlCenter
$a1.VerticalAlignment = -4160    # xlTop
$a1.Borders.LineStyle = 1        # xlContinuous
$a1.Borders.Weight = 2           # xlThin

# Copy the exact same formatting onto B1 (reuses the same style record
# instead of spawning extra ones for an equivalent but distinct cellXf).
$a1.Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
